$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # RoundTrip
$ws2 = $wb.Worksheets.Item(2)   # HotelSearch
$ws3 = $wb.Worksheets.Item(3)   # Sheet3

# --- Sheet3: add the Yes/No source list used by the dropdown ---
$ws3.Range("A16").Value = "Yes"
$ws3.Range("A17").Value = "No"

# --- Sheet1 (RoundTrip): rework the Yes/No helper list + dropdown values ---
# Shrink the hidden helper column A width
$ws1.Columns.Item(1).ColumnWidth = 11.022135416666666

# New helper list lives in A3:A4 instead of A8:A9
$ws1.Range("A8").ClearContents()
$ws1.Range("A9").ClearContents()
$ws1.Range("A3").Value = "Yes"
$ws1.Range("A4").Value = "No"

# Update the I-column dropdown values to match
$ws1.Range("I2").Value = "No"
$ws1.Range("I3").Value = "Yes"
$ws1.Range("I4").Value = "Yes"

# Re-point the dropdown validation at the new helper list + extend to I1
$i_range = $ws1.Range("I1:I1048576")
$i_range.Validation.Delete()
$i_range.Validation.Add(3, 1, 1, "=`$A`$3:`$A`$4")

# --- Sheet2 (HotelSearch): shrink default/column widths ---
$ws2.Columns.Item(1).ColumnWidth = -0.8333333333333334
$ws2.Range("B1:F1").EntireColumn.ColumnWidth = 16.166666666666668

# --- Selections / active sheet (applied last, in tab order) ---
$ws3.Select() | Out-Null
$ws3.Range("C25").Select() | Out-Null

$ws2.Select() | Out-Null
$ws2.Range("C12").Select() | Out-Null

$ws1.Select() | Out-Null
$ws1.Range("E14").Select() | Out-Null
